$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Email column values (D2:D5) were removed from the sheet.
$ws.Range("D2:D5").ClearContents()

# Leave the sheet with D2:D5 selected, matching the state after the delete.
$ws.Range("D2:D5").Select()
